# This workbook tracks "Avverkningsanmälningar" (logging notifications).
# The source data was refreshed: the "Förändrad" (changed) column for every
# data row (2-511) is bumped from 2023-09-20 (45189) to 2023-09-21 (45190).
# Additionally, the last four rows (508-511), which all share the same
# "Datum" (45188) and therefore were tied, got re-ordered (re-sorted by
# Area ascending) as part of that refresh, so their "Beteckning" (A) and
# "Area (ha)" (G) values need to be updated in place to match the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" date for every data row (2 through 511) from
# 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C511").Value = 45190

# The final four rows were re-sorted (ascending by Area). Update the
# "Beteckning" (column A) and "Area (ha)" (column G) values in place so the
# rows end up holding the post-sort data.
$ws.Cells.Item(508, 1).Value = "A 44090-2023"
$ws.Cells.Item(508, 7).Value = 0.5

$ws.Cells.Item(509, 1).Value = "A 44087-2023"
$ws.Cells.Item(509, 7).Value = 1.1

$ws.Cells.Item(510, 1).Value = "A 44086-2023"
$ws.Cells.Item(510, 7).Value = 2.1

$ws.Cells.Item(511, 1).Value = "A 44133-2023"
$ws.Cells.Item(511, 7).Value = 5.1
